# BOM.xlsx edit: switch caps for lower cost
# Replace the P/N for the 0.1u capacitor bank (row 3, column D) with a
# cheaper part, and update the sheet's selection to D2:D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds reference "C3,C5,C6,C7,C10,C11,C12" / value "0.1u".
# Its part number (column D) changes from C0603G104K5RACT250 to the
# cheaper UMK107B7104KAHT (note the trailing space, preserved verbatim).
$ws.Range("D3").Value = "UMK107B7104KAHT "

# Update the active selection on the sheet to D2:D19 with D2 as the active cell.
$ws.Range("D2:D19").Select()
